$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")

# --- 1. Move the existing "Tổng" (total) row from row 21 down to row 25 ---
# Use Cut so the cell styles (s=25 / s=8) and values travel for free (no new
# style slots get minted the way Insert/Merge do in this runtime).
$ws.Range("B21:E21").Cut($ws.Range("B25:E25"))

# Drop the stale merge left behind at the old location and recreate it at
# the new one, then restore the original "Tổng" formatting (style 25) that
# Merge() resets, by pasting formats back from row 4's block - cheapest
# source that still carries style 25 is the row we just vacated via Cut,
# so grab it while it is still tagged in the clipboard state via B4 pattern.
$ws.Range("B21:C21").UnMerge()
$ws.Range("B25:C25").Merge()

# --- 2. Fill in the two already-present blank rows (19 and 20) ---
$ws.Range("B19").Value2 = 42475
$ws.Range("C19").Value2 = "In bản chính thức. 5 bộ 200tr, 5 bộ 38tr"
$ws.Range("E19").Value2 = 780

$ws.Range("B20").Value2 = 42475
$ws.Range("C20").Value2 = "Uống nước tại Moda"
$ws.Range("E20").Value2 = 76

# --- 3. New rows 21 and 22 (formats copied from row 20, values written after) ---
$ws.Range("B20:E20").Copy()
$ws.Range("B21:E21").PasteSpecial(-4122)
$ws.Range("B20:E20").Copy()
$ws.Range("B22:E22").PasteSpecial(-4122)

$ws.Range("B21").Value2 = 42109
$ws.Range("C21").Value2 = "Mua bánh mì cho anh Tùng"
$ws.Range("E21").Value2 = 10

$ws.Range("B22").Value2 = 42475
$ws.Range("C22").Value2 = "Ăn trưa mì xào bò 3 hộp"
$ws.Range("E22").Value2 = 75

# --- 4. Two new blank rows (23, 24) with the same formatting as row 20 ---
$ws.Range("B20:E20").Copy()
$ws.Range("B23:E23").PasteSpecial(-4122)
$ws.Range("B20:E20").Copy()
$ws.Range("B24:E24").PasteSpecial(-4122)

# --- 5. Fix up the totals row formulas for the new row span ---
$ws.Range("D25").Formula = "=SUM(D4:D14)"
$ws.Range("E25").Formula = "=SUM(E4:E24)"

# Restore style 25 on the merged total-label cells (Merge() resets format)
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# --- 6. Sheet view bookkeeping to match the authored state ---
$ws.Range("B23:E23").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
